$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text columns (Coin name / Link URL) ---
$ws.Range('B24').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C24').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('B30').Value = 'LidoDAOToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('B38').Value = 'Frax'
$ws.Range('C38').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('B48').Value = 'PEPE'
$ws.Range('C48').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'

# --- Numeric-looking text columns (Price / Volume) ---
# Force text number format first so Excel does not coerce these
# digit/percent strings into actual numbers.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.990.73'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.38%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.913.23'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.71%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '324.72'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.31%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4595'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.75%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07724'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9822'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '22.14'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.54%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.919.57'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.952'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.93%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.669'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.07031'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.26%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '83.90'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -3.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009468'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -3.71%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.70'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.80%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '28.966.14'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.65%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.322'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -2.75%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.26%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.130.86'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.98%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.088'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '158.43'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.72%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.03'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.72%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.673'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.47%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '117.45'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.94%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.855'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.40%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09305'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.82%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.8695'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.49%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.081'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.87%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.253'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -4.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.026'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -3.17%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05729'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.50%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.155'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.40%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.001'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.09%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02044'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.95%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5509'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.68%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.403'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.71%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1754'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.54%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.847'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +4.15%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '9.339'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.35%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5179'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -2.11%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '11.26'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.91%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.06871'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.000002623'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -9.60%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.055'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -3.44%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.785'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.76%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '110.63'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.99%  '

# --- Strip the explicit "@" text format back off so the cells keep
#     their original default styling (no style index), matching
#     the look of the untouched cells around them. ---
$deCells = @('D2','E2','D3','E3','D4','E4','D5','E5','E6','D7','E7','E8','D9','E9','D10','E10','D11','E11','D12','E12','D13','E13','D14','E14','D15','E15','E16','D17','E17','D18','E18','D19','E19','E20','D21','E21','D22','E22','E23','D24','E24','D25','E25','D26','E26','D27','E27','D28','E28','D29','E29','D30','E30','D31','E31','D32','E32','D33','E33','D34','E34','D35','E35','D36','E36','D37','E37','D38','E38','D39','E39','D40','E40','D41','E41','D42','E42','D43','E43','D44','E44','D45','E45','D46','E46','D47','E47','D48','E48','D49','E49','D50','E50','D51','E51')
foreach ($ref in $deCells) {
    $ws.Range($ref).ClearFormats()
}